# localization-status.xlsx : "Generate Report for Archive"
#
# The status of the single tracked file flips from "Ready for handoff" to
# "In Translation" everywhere it is shown (Overview!E2:F2, zh-cn!C2,
# de-de!C2), and the now-narrower text means the Status columns can be
# drawn tighter, so their widths shrink to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update every cell that shows the old status so the shared string is
# replaced cleanly (no cell is left pointing at "Ready for handoff").
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Narrow the Status columns to fit the shorter text.
$overview.Columns(5).ColumnWidth = 12.58   # column E
$overview.Columns(6).ColumnWidth = 12.58   # column F
$zhcn.Columns(3).ColumnWidth = 12.58       # column C
$dede.Columns(3).ColumnWidth = 12.58       # column C
